$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the leftover "_GoBack" bookmark (Word inserts/keeps this marker
#    around the last edit location; it was cleaned up in this revision).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Add a new line: "I add a line F. from perth", surrounded by a blank
#    paragraph before and after, inserted right before the two trailing
#    empty paragraphs at the end of the document (and before the sectPr).
# ---------------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$anchor = $d.Paragraphs.Item($paraCount - 1)
$insertionPoint = $d.Range($anchor.Range.Start, $anchor.Range.Start)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newParagraphsXml = (
    "<w:p $wNs/>" +
    "<w:p $wNs>" +
        "<w:r><w:t xml:space=`"preserve`">I add a line F. from </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/>" +
        "<w:r><w:t>perth</w:t></w:r>" +
        "<w:proofErr w:type=`"spellEnd`"/>" +
    "</w:p>" +
    "<w:p $wNs/>"
)

$null = $insertionPoint.InsertXML($newParagraphsXml)
